# إضافة حدث جديد في Card24 by admin at 2025-12-08 07:01:36
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")

# --- Fill the previously-blank placeholder cells on row 17 with "nan" ---
$ws.Range("B17").Value = "nan"
$ws.Range("C17").Value = "nan"
$ws.Range("D17").Value = "nan"
$ws.Range("E17").Value = "nan"
$ws.Range("F17").Value = "nan"
$ws.Range("G17").Value = "nan"
$ws.Range("H17").Value = "nan"
$ws.Range("I17").Value = "nan"
$ws.Range("J17").Value = "nan"
$ws.Range("K17").Value = "nan"
$ws.Range("P17").Value = "nan"

# --- Add the new service event as row 18 ---
# Column A holds a numeric-looking card id that must stay a text value
# (matches how "24" is stored as text elsewhere in this sheet).
$ws.Range("A18").NumberFormat = "@"
$ws.Range("A18").Value = "24"
$ws.Range("A18").ClearFormats()

$ws.Range("L18").Value = "6\10\2025"
$ws.Range("M18").Value = "870 t"
$ws.Range("N18").Value = "تم تغيير زيت الجروبوكس الفلتس وتغير الجرائد الخلفيه (1_5_8)"
$ws.Range("O18").Value = "تم العمل"

Write-Host "Row 17 filled and row 18 added to Card24"
